$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first data table
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 3526
$wsExpo.Range("F5").Value  = 2192
$wsExpo.Range("F6").Value  = 428
$wsExpo.Range("F7").Value  = 169
$wsExpo.Range("F8").Value  = 68
$wsExpo.Range("F9").Value  = 57
$wsExpo.Range("F10").Value = 1301
$wsExpo.Range("F12").Value = 1772

# Sheet "全部类型" (All types) - combined data table
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 3526
$wsAll.Range("F5").Value  = 2192
$wsAll.Range("F6").Value  = 428
$wsAll.Range("F8").Value  = 169
$wsAll.Range("F9").Value  = 68
$wsAll.Range("F10").Value = 57
$wsAll.Range("F13").Value = 1301
$wsAll.Range("F15").Value = 1772
